$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.624282666666667
$ws.Range("H2").Value = 4.872847999999999
$ws.Range("I2").Value = 0.08561172663893989
$ws.Range("J2").Value = 0.08561172663893987
$ws.Range("M2").Value = 9.172748
$ws.Range("N2").Value = 27.518244
$ws.Range("O2").Value = 0.01445826353606064
$ws.Range("P2").Value = 0.01445826353606064
$ws.Range("Q2").Value = 14.89913558210133
$ws.Range("R2").Value = 134.092220238912
$ws.Range("S2").Value = 0.001237796905522976
$ws.Range("T2").Value = 0.001237796905522975
$ws.Range("G3").Value = 1.624282666666667
$ws.Range("H3").Value = 4.872847999999999
$ws.Range("I3").Value = 0.08561172663893989
$ws.Range("J3").Value = 0.08561172663893987
$ws.Range("O3").Value = 0.2254554169720557
$ws.Range("P3").Value = 0.2254554169720557
$ws.Range("Q3").Value = 232.3301699964089
$ws.Range("R3").Value = 2090.97152996768
$ws.Range("S3").Value = 0.01930162752707984
$ws.Range("T3").Value = 0.01930162752707984
$ws.Range("G4").Value = 1.624282666666667
$ws.Range("H4").Value = 4.872847999999999
$ws.Range("I4").Value = 0.08561172663893989
$ws.Range("J4").Value = 0.08561172663893987
$ws.Range("M4").Value = 169.2367096666667
$ws.Range("N4").Value = 507.7101290000001
$ws.Range("O4").Value = 0.2667541884216647
$ws.Range("P4").Value = 0.2667541884216647
$ws.Range("Q4").Value = 274.8882540752658
$ws.Range("R4").Value = 2473.994286677392
$ws.Range("S4").Value = 0.02283728665894782
$ws.Range("T4").Value = 0.02283728665894781
$ws.Range("G5").Value = 1.624282666666667
$ws.Range("H5").Value = 4.872847999999999
$ws.Range("I5").Value = 0.08561172663893989
$ws.Range("J5").Value = 0.08561172663893987
$ws.Range("M5").Value = 16.15031566666667
$ws.Range("N5").Value = 48.450947
$ws.Range("O5").Value = 0.02545644119943506
$ws.Range("P5").Value = 0.02545644119943505
$ws.Range("Q5").Value = 26.23267779856178
$ws.Range("R5").Value = 236.094100187056
$ws.Range("S5").Value = 0.002179369885166281
$ws.Range("T5").Value = 0.00217936988516628
$ws.Range("G6").Value = 1.624282666666667
$ws.Range("H6").Value = 4.872847999999999
$ws.Range("I6").Value = 0.08561172663893989
$ws.Range("J6").Value = 0.08561172663893987
$ws.Range("M6").Value = 54.744643
$ws.Range("N6").Value = 164.233929
$ws.Range("O6").Value = 0.08628956945961638
$ws.Range("P6").Value = 0.08628956945961638
$ws.Range("Q6").Value = 88.92077471775467
$ws.Range("R6").Value = 800.286972459792
$ws.Range("S6").Value = 0.007387399032368494
$ws.Range("T6").Value = 0.007387399032368493
$ws.Range("G7").Value = 1.624282666666667
$ws.Range("H7").Value = 4.872847999999999
$ws.Range("I7").Value = 0.08561172663893989
$ws.Range("J7").Value = 0.08561172663893987
$ws.Range("M7").Value = 242.0894676666667
$ws.Range("N7").Value = 726.268403
$ws.Range("O7").Value = 0.3815861204111676
$ws.Range("P7").Value = 0.3815861204111676
$ws.Range("Q7").Value = 393.2217261135271
$ws.Range("R7").Value = 3538.995535021744
$ws.Range("S7").Value = 0.03266824662985449
$ws.Range("T7").Value = 0.03266824662985448
$ws.Range("I8").Value = 0.1920894545885022
$ws.Range("J8").Value = 0.1920894545885022
$ws.Range("M8").Value = 9.172748
$ws.Range("N8").Value = 27.518244
$ws.Range("O8").Value = 0.01445826353606064
$ws.Range("P8").Value = 0.01445826353606064
$ws.Range("Q8").Value = 33.429612275852
$ws.Range("R8").Value = 300.866510482668
$ws.Range("S8").Value = 0.002777279956938716
$ws.Range("T8").Value = 0.002777279956938716
$ws.Range("I9").Value = 0.1920894545885022
$ws.Range("J9").Value = 0.1920894545885022
$ws.Range("O9").Value = 0.2254554169720557
$ws.Range("P9").Value = 0.2254554169720557
$ws.Range("S9").Value = 0.04330760808018552
$ws.Range("T9").Value = 0.04330760808018551
$ws.Range("I10").Value = 0.1920894545885022
$ws.Range("J10").Value = 0.1920894545885022
$ws.Range("M10").Value = 169.2367096666667
$ws.Range("N10").Value = 507.7101290000001
$ws.Range("O10").Value = 0.2667541884216647
$ws.Range("P10").Value = 0.2667541884216647
$ws.Range("Q10").Value = 616.7745573079737
$ws.Range("R10").Value = 5550.971015771765
$ws.Range("S10").Value = 0.0512406665631161
$ws.Range("T10").Value = 0.0512406665631161
$ws.Range("I11").Value = 0.1920894545885022
$ws.Range("J11").Value = 0.1920894545885022
$ws.Range("M11").Value = 16.15031566666667
$ws.Range("N11").Value = 48.450947
$ws.Range("O11").Value = 0.02545644119943506
$ws.Range("P11").Value = 0.02545644119943505
$ws.Range("Q11").Value = 58.85900178106768
$ws.Range("R11").Value = 529.731016029609
$ws.Range("S11").Value = 0.004889913905763756
$ws.Range("T11").Value = 0.004889913905763755
$ws.Range("I12").Value = 0.1920894545885022
$ws.Range("J12").Value = 0.1920894545885022
$ws.Range("M12").Value = 54.744643
$ws.Range("N12").Value = 164.233929
$ws.Range("O12").Value = 0.08628956945961638
$ws.Range("P12").Value = 0.08628956945961638
$ws.Range("Q12").Value = 199.514059436707
$ws.Range("R12").Value = 1795.626534930364
$ws.Range("S12").Value = 0.01657531633417438
$ws.Range("T12").Value = 0.01657531633417438
$ws.Range("I13").Value = 0.1920894545885022
$ws.Range("J13").Value = 0.1920894545885022
$ws.Range("M13").Value = 242.0894676666667
$ws.Range("N13").Value = 726.268403
$ws.Range("O13").Value = 0.3815861204111676
$ws.Range("P13").Value = 0.3815861204111676
$ws.Range("Q13").Value = 882.2827183483158
$ws.Range("R13").Value = 7940.544465134843
$ws.Range("S13").Value = 0.0732986697483237
$ws.Range("T13").Value = 0.0732986697483237
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.741573333333333
$ws.Range("H14").Value = 8.22472
$ws.Range("I14").Value = 0.1445012198865677
$ws.Range("J14").Value = 0.1445012198865677
$ws.Range("M14").Value = 9.172748
$ws.Range("N14").Value = 27.518244
$ws.Range("O14").Value = 0.01445826353606064
$ws.Range("P14").Value = 0.01445826353606064
$ws.Range("Q14").Value = 25.14776131018667
$ws.Range("R14").Value = 226.32985179168
$ws.Range("S14").Value = 0.002089236718402242
$ws.Range("T14").Value = 0.002089236718402242
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.741573333333333
$ws.Range("H15").Value = 8.22472
$ws.Range("I15").Value = 0.1445012198865677
$ws.Range("J15").Value = 0.1445012198865677
$ws.Range("O15").Value = 0.2254554169720557
$ws.Range("P15").Value = 0.2254554169720557
$ws.Range("Q15").Value = 392.1424587372445
$ws.Range("R15").Value = 3529.2821286352
$ws.Range("S15").Value = 0.03257858278249683
$ws.Range("T15").Value = 0.03257858278249683
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.741573333333333
$ws.Range("H16").Value = 8.22472
$ws.Range("I16").Value = 0.1445012198865677
$ws.Range("J16").Value = 0.1445012198865677
$ws.Range("M16").Value = 169.2367096666667
$ws.Range("N16").Value = 507.7101290000001
$ws.Range("O16").Value = 0.2667541884216647
$ws.Range("P16").Value = 0.2667541884216647
$ws.Range("Q16").Value = 463.9748502432089
$ws.Range("R16").Value = 4175.77365218888
$ws.Range("S16").Value = 0.03854630563678187
$ws.Range("T16").Value = 0.03854630563678187
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.741573333333333
$ws.Range("H17").Value = 8.22472
$ws.Range("I17").Value = 0.1445012198865677
$ws.Range("J17").Value = 0.1445012198865677
$ws.Range("M17").Value = 16.15031566666667
$ws.Range("N17").Value = 48.450947
$ws.Range("O17").Value = 0.02545644119943506
$ws.Range("P17").Value = 0.02545644119943505
$ws.Range("Q17").Value = 44.27727475664889
$ws.Range("R17").Value = 398.49547280984
$ws.Range("S17").Value = 0.003678486807289046
$ws.Range("T17").Value = 0.003678486807289045
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 2.741573333333333
$ws.Range("H18").Value = 8.22472
$ws.Range("I18").Value = 0.1445012198865677
$ws.Range("J18").Value = 0.1445012198865677
$ws.Range("M18").Value = 54.744643
$ws.Range("N18").Value = 164.233929
$ws.Range("O18").Value = 0.08628956945961638
$ws.Range("P18").Value = 0.08628956945961638
$ws.Range("Q18").Value = 150.0864533916533
$ws.Range("R18").Value = 1350.77808052488
$ws.Range("S18").Value = 0.01246894805040128
$ws.Range("T18").Value = 0.01246894805040128
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 2.741573333333333
$ws.Range("H19").Value = 8.22472
$ws.Range("I19").Value = 0.1445012198865677
$ws.Range("J19").Value = 0.1445012198865677
$ws.Range("M19").Value = 242.0894676666667
$ws.Range("N19").Value = 726.268403
$ws.Range("O19").Value = 0.3815861204111676
$ws.Range("P19").Value = 0.3815861204111676
$ws.Range("Q19").Value = 663.7060288357956
$ws.Range("R19").Value = 5973.35425952216
$ws.Range("S19").Value = 0.05513965989119643
$ws.Range("T19").Value = 0.05513965989119643
$ws.Range("G20").Value = 2.34427
$ws.Range("H20").Value = 7.03281
$ws.Range("I20").Value = 0.1235603916279767
$ws.Range("J20").Value = 0.1235603916279767
$ws.Range("M20").Value = 9.172748
$ws.Range("N20").Value = 27.518244
$ws.Range("O20").Value = 0.01445826353606064
$ws.Range("P20").Value = 0.01445826353606064
$ws.Range("Q20").Value = 21.50339795396
$ws.Range("R20").Value = 193.53058158564
$ws.Range("S20").Value = 0.001786468704776147
$ws.Range("T20").Value = 0.001786468704776147
$ws.Range("G21").Value = 2.34427
$ws.Range("H21").Value = 7.03281
$ws.Range("I21").Value = 0.1235603916279767
$ws.Range("J21").Value = 0.1235603916279767
$ws.Range("O21").Value = 0.2254554169720557
$ws.Range("P21").Value = 0.2254554169720557
$ws.Range("Q21").Value = 335.3139566127334
$ws.Range("R21").Value = 3017.8256095146
$ws.Range("S21").Value = 0.02785735961571598
$ws.Range("T21").Value = 0.02785735961571598
$ws.Range("G22").Value = 2.34427
$ws.Range("H22").Value = 7.03281
$ws.Range("I22").Value = 0.1235603916279767
$ws.Range("J22").Value = 0.1235603916279767
$ws.Range("M22").Value = 169.2367096666667
$ws.Range("N22").Value = 507.7101290000001
$ws.Range("O22").Value = 0.2667541884216647
$ws.Range("P22").Value = 0.2667541884216647
$ws.Range("Q22").Value = 396.7365413702767
$ws.Range("R22").Value = 3570.62887233249
$ws.Range("S22").Value = 0.03296025198978396
$ws.Range("T22").Value = 0.03296025198978396
$ws.Range("G23").Value = 2.34427
$ws.Range("H23").Value = 7.03281
$ws.Range("I23").Value = 0.1235603916279767
$ws.Range("J23").Value = 0.1235603916279767
$ws.Range("M23").Value = 16.15031566666667
$ws.Range("N23").Value = 48.450947
$ws.Range("O23").Value = 0.02545644119943506
$ws.Range("P23").Value = 0.02545644119943505
$ws.Range("Q23").Value = 37.86070050789667
$ws.Range("R23").Value = 340.74630457107
$ws.Range("S23").Value = 0.003145407844056755
$ws.Range("T23").Value = 0.003145407844056755
$ws.Range("G24").Value = 2.34427
$ws.Range("H24").Value = 7.03281
$ws.Range("I24").Value = 0.1235603916279767
$ws.Range("J24").Value = 0.1235603916279767
$ws.Range("M24").Value = 54.744643
$ws.Range("N24").Value = 164.233929
$ws.Range("O24").Value = 0.08628956945961638
$ws.Range("P24").Value = 0.08628956945961638
$ws.Range("Q24").Value = 128.33622424561
$ws.Range("R24").Value = 1155.02601821049
$ws.Range("S24").Value = 0.01066197299583969
$ws.Range("T24").Value = 0.01066197299583969
$ws.Range("G25").Value = 2.34427
$ws.Range("H25").Value = 7.03281
$ws.Range("I25").Value = 0.1235603916279767
$ws.Range("J25").Value = 0.1235603916279767
$ws.Range("M25").Value = 242.0894676666667
$ws.Range("N25").Value = 726.268403
$ws.Range("O25").Value = 0.3815861204111676
$ws.Range("P25").Value = 0.3815861204111676
$ws.Range("Q25").Value = 567.5230763669366
$ws.Range("R25").Value = 5107.70768730243
$ws.Range("S25").Value = 0.04714893047780413
$ws.Range("T25").Value = 0.04714893047780413
$ws.Range("G26").Value = 4.282534999999999
$ws.Range("H26").Value = 12.847605
$ws.Range("I26").Value = 0.225721312715906
$ws.Range("J26").Value = 0.225721312715906
$ws.Range("M26").Value = 9.172748
$ws.Range("N26").Value = 27.518244
$ws.Range("O26").Value = 0.01445826353606064
$ws.Range("P26").Value = 0.01445826353606064
$ws.Range("Q26").Value = 39.28261435618
$ws.Range("R26").Value = 353.54352920562
$ws.Range("S26").Value = 0.003263538224952124
$ws.Range("T26").Value = 0.003263538224952124
$ws.Range("G27").Value = 4.282534999999999
$ws.Range("H27").Value = 12.847605
$ws.Range("I27").Value = 0.225721312715906
$ws.Range("J27").Value = 0.225721312715906
$ws.Range("O27").Value = 0.2254554169720557
$ws.Range("P27").Value = 0.2254554169720557
$ws.Range("Q27").Value = 612.5547633943667
$ws.Range("R27").Value = 5512.992870549299
$ws.Range("S27").Value = 0.05089009267784437
$ws.Range("T27").Value = 0.05089009267784436
$ws.Range("G28").Value = 4.282534999999999
$ws.Range("H28").Value = 12.847605
$ws.Range("I28").Value = 0.225721312715906
$ws.Range("J28").Value = 0.225721312715906
$ws.Range("M28").Value = 169.2367096666667
$ws.Range("N28").Value = 507.7101290000001
$ws.Range("O28").Value = 0.2667541884216647
$ws.Range("P28").Value = 0.2667541884216647
$ws.Range("Q28").Value = 724.7621324323383
$ws.Range("R28").Value = 6522.859191891044
$ws.Range("S28").Value = 0.06021210558300428
$ws.Range("T28").Value = 0.06021210558300428
$ws.Range("G29").Value = 4.282534999999999
$ws.Range("H29").Value = 12.847605
$ws.Range("I29").Value = 0.225721312715906
$ws.Range("J29").Value = 0.225721312715906
$ws.Range("M29").Value = 16.15031566666667
$ws.Range("N29").Value = 48.450947
$ws.Range("O29").Value = 0.02545644119943506
$ws.Range("P29").Value = 0.02545644119943505
$ws.Range("Q29").Value = 69.16429210354832
$ws.Range("R29").Value = 622.4786289319348
$ws.Range("S29").Value = 0.005746061324611753
$ws.Range("T29").Value = 0.005746061324611753
$ws.Range("G30").Value = 4.282534999999999
$ws.Range("H30").Value = 12.847605
$ws.Range("I30").Value = 0.225721312715906
$ws.Range("J30").Value = 0.225721312715906
$ws.Range("M30").Value = 54.744643
$ws.Range("N30").Value = 164.233929
$ws.Range("O30").Value = 0.08628956945961638
$ws.Range("P30").Value = 0.08628956945961638
$ws.Range("Q30").Value = 234.445849710005
$ws.Range("R30").Value = 2110.012647390045
$ws.Range("S30").Value = 0.01947739489211496
$ws.Range("T30").Value = 0.01947739489211496
$ws.Range("G31").Value = 4.282534999999999
$ws.Range("H31").Value = 12.847605
$ws.Range("I31").Value = 0.225721312715906
$ws.Range("J31").Value = 0.225721312715906
$ws.Range("M31").Value = 242.0894676666667
$ws.Range("N31").Value = 726.268403
$ws.Range("O31").Value = 0.3815861204111676
$ws.Range("P31").Value = 0.3815861204111676
$ws.Range("Q31").Value = 1036.756618413868
$ws.Range("R31").Value = 9330.809565724814
$ws.Range("S31").Value = 0.08613212001337853
$ws.Range("T31").Value = 0.08613212001337853
$ws.Range("E32").Value = 3
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = 4.335555666666667
$ws.Range("H32").Value = 13.006667
$ws.Range("I32").Value = 0.2285158945421077
$ws.Range("J32").Value = 0.2285158945421077
$ws.Range("M32").Value = 9.172748
$ws.Range("N32").Value = 27.518244
$ws.Range("O32").Value = 0.01445826353606064
$ws.Range("P32").Value = 0.01445826353606064
$ws.Range("Q32").Value = 39.76895957030534
$ws.Range("R32").Value = 357.920636132748
$ws.Range("S32").Value = 0.003303943025468434
$ws.Range("T32").Value = 0.003303943025468433
$ws.Range("E33").Value = 3
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 4.335555666666667
$ws.Range("H33").Value = 13.006667
$ws.Range("I33").Value = 0.2285158945421077
$ws.Range("J33").Value = 0.2285158945421077
$ws.Range("O33").Value = 0.2254554169720557
$ws.Range("P33").Value = 0.2254554169720557
$ws.Range("Q33").Value = 620.1386037891357
$ws.Range("R33").Value = 5581.247434102221
$ws.Range("S33").Value = 0.05152014628873321
$ws.Range("T33").Value = 0.05152014628873319
$ws.Range("E34").Value = 3
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 4.335555666666667
$ws.Range("H34").Value = 13.006667
$ws.Range("I34").Value = 0.2285158945421077
$ws.Range("J34").Value = 0.2285158945421077
$ws.Range("M34").Value = 169.2367096666667
$ws.Range("N34").Value = 507.7101290000001
$ws.Range("O34").Value = 0.2667541884216647
$ws.Range("P34").Value = 0.2667541884216647
$ws.Range("Q34").Value = 733.7351756033382
$ws.Range("R34").Value = 6603.616580430044
$ws.Range("S34").Value = 0.06095757199003065
$ws.Range("T34").Value = 0.06095757199003064
$ws.Range("E35").Value = 3
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 4.335555666666667
$ws.Range("H35").Value = 13.006667
$ws.Range("I35").Value = 0.2285158945421077
$ws.Range("J35").Value = 0.2285158945421077
$ws.Range("M35").Value = 16.15031566666667
$ws.Range("N35").Value = 48.450947
$ws.Range("O35").Value = 0.02545644119943506
$ws.Range("P35").Value = 0.02545644119943505
$ws.Range("Q35").Value = 70.02059260707212
$ws.Range("R35").Value = 630.185333463649
$ws.Range("S35").Value = 0.005817201432547467
$ws.Range("T35").Value = 0.005817201432547466
$ws.Range("E36").Value = 3
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 4.335555666666667
$ws.Range("H36").Value = 13.006667
$ws.Range("I36").Value = 0.2285158945421077
$ws.Range("J36").Value = 0.2285158945421077
$ws.Range("M36").Value = 54.744643
$ws.Range("N36").Value = 164.233929
$ws.Range("O36").Value = 0.08628956945961638
$ws.Range("P36").Value = 0.08628956945961638
$ws.Range("Q36").Value = 237.3484471782937
$ws.Range("R36").Value = 2136.136024604643
$ws.Range("S36").Value = 0.01971853815471758
$ws.Range("T36").Value = 0.01971853815471757
$ws.Range("E37").Value = 3
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 4.335555666666667
$ws.Range("H37").Value = 13.006667
$ws.Range("I37").Value = 0.2285158945421077
$ws.Range("J37").Value = 0.2285158945421077
$ws.Range("M37").Value = 242.0894676666667
$ws.Range("N37").Value = 726.268403
$ws.Range("O37").Value = 0.3815861204111676
$ws.Range("P37").Value = 0.3815861204111676
$ws.Range("Q37").Value = 1049.592363382534
$ws.Range("R37").Value = 9446.3312704428
$ws.Range("S37").Value = 0.0871984936506104
$ws.Range("T37").Value = 0.08719849365061039
